$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "3433 HammondCare Caulfield Village Aged Care" cluster row entirely.
$ws.Rows("2:2").Delete()

# Fix typo in remaining cluster name (missing space before "Department").
$ws.Cells.Item(12, 1).Value = "Western Health Sunshine Hospital Emergency Department St Albans"

# Refresh the active case counts for the remaining clusters.
$ws.Cells.Item(2, 2).Value = 38
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(4, 2).Value = 20
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(6, 2).Value = 20
$ws.Cells.Item(7, 2).Value = 38
$ws.Cells.Item(8, 2).Value = 13
$ws.Cells.Item(9, 2).Value = 14
$ws.Cells.Item(10, 2).Value = 13
$ws.Cells.Item(11, 2).Value = 16
$ws.Cells.Item(12, 2).Value = 10
